# Apply updated cryptocurrency data (prices / 1h volume %) scraped on
# Thu Jul 20 13:53:14 UTC 2023, matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every touched cell so numeric-looking strings
# (e.g. "30.287.23", "0.000007869", "8.250") are kept as literal text,
# matching the inlineStr cells in the original workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.287.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8072'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.76'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3265'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.04'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07276'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7896'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08095'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.917.04'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.418'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.14'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.272.19'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.083'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.63'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007869'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.250'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +20.51%  '
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.168.63'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1638'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +15.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.498'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.14'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.05'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.164'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.395'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.415'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05732'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.155'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.301'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7505'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.007'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.732'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01963'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.823'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4548'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.013'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8552'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.936'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.66'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.035.32'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.02'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.641'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.098'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +11.02%  '
